# Regenerate save_data to use K (strikeouts) instead of Strike# values.
# Column G holds the "K" stat; update rows 2-12 with the recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 7
    3  = 6
    4  = 3
    5  = 5
    6  = 9
    7  = 4
    8  = 8
    9  = 7
    10 = 2
    11 = 3
    12 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
